# Update the date line and the division problems in the table.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Date heading
Replace-Text "2025-05-12 Monday" "2025-05-13 Tuesday"

# The table has a single <w:tbl> with 20 rows; rows 1, 5, 9, 13, 17 (1-based)
# hold the 5 division problems each. Address cells directly by row/column so
# the two identical "81÷4=" occurrences (row1/col2 and row13/col4) resolve
# to their distinct replacements.
$table = $d.Tables.Item(1)

$values = @{
    1  = @("30÷3=", "94÷4=", "13÷4=", "95÷4=", "79÷7=")
    5  = @("70÷8=", "14÷2=", "17÷3=", "47÷7=", "26÷6=")
    9  = @("56÷8=", "60÷6=", "78÷7=", "83÷4=", "75÷7=")
    13 = @("88÷5=", "63÷6=", "62÷4=", "18÷9=", "27÷5=")
    17 = @("84÷8=", "13÷5=", "25÷2=", "56÷6=", "30÷5=")
}

foreach ($rowIndex in $values.Keys) {
    $cols = $values[$rowIndex]
    for ($c = 1; $c -le $cols.Length; $c++) {
        $cell = $table.Cell($rowIndex, $c)
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1) | Out-Null   # drop the end-of-cell marker
        $cellRange.Text = $cols[$c - 1]
    }
}
